# Home page: changed title + added description.
# Append two blank paragraphs (language fr-FR) right after the
# "Brève description des fonctionnalités" Heading 1, i.e. right before
# the final section break, mirroring pressing Enter twice at the end of
# that heading (new paragraphs pick up the Normal "next style" with no
# visible text, only the paragraph-mark language formatting).

$d = $word.ActiveDocument

# Locate the heading unambiguously via its bookmark (the same heading
# also appears as a TOC entry earlier in the document, so a plain text
# search would be ambiguous).
$bm = $d.Bookmarks("_Toc388547958")
$headingPara = $d.Paragraphs.Last
if ($bm.Range.Start -lt $headingPara.Range.Start) {
    # Fallback safety net, should not normally trigger.
    $headingPara = $bm.Range.Paragraphs.Last
}

$insertPoint = $d.Range($headingPara.Range.End, $headingPara.Range.End)

$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr></w:p>'
$insertPoint.InsertXML($blankParaXml + $blankParaXml)
